$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 275.2143
$ws.Range("I12").Value = 275.2143
$ws.Range("K12").Value = 275.2143
$ws.Range("M12").Value = -105.2143
$ws.Range("H74").Value = 5912.5
$ws.Range("I74").Value = 4743
$ws.Range("K74").Value = 4743
$ws.Range("M74").Value = -3807
$ws.Range("H77").Value = 5912.5
$ws.Range("I77").Value = 4743
$ws.Range("K77").Value = 23747.5
$ws.Range("M77").Value = -19035
$ws.Range("H100").Value = 89295.88
$ws.Range("I100").Value = 125872.375
$ws.Range("J100").Value = 56783.445
$ws.Range("K100").Value = 125872.375
$ws.Range("L100").Value = 56783.445
$ws.Range("M100").Value = -125331.375
$ws.Range("N100").Value = -57865.445
$ws.Range("H113").Value = 4921.643
$ws.Range("I113").Value = 5167
$ws.Range("J113").Value = 3449.5
$ws.Range("K113").Value = 5167
$ws.Range("L113").Value = 3449.5
$ws.Range("M113").Value = -1913
$ws.Range("N113").Value = -9957.5
$ws.Range("H132").Value = 2556.5264
$ws.Range("I132").Value = 2447.018
$ws.Range("K132").Value = 7341.054
$ws.Range("M132").Value = -4811.054
$ws.Range("H135").Value = 5176.484
$ws.Range("I135").Value = 5859.3335
$ws.Range("K135").Value = 52734.0015
$ws.Range("M135").Value = -50199.0015
$ws.Range("H137").Value = 8638.725
$ws.Range("I137").Value = 13118.608
$ws.Range("K137").Value = 39355.824
$ws.Range("M137").Value = -36805.824
$ws.Range("H138").Value = 4430.909
$ws.Range("I138").Value = 2232.3333
$ws.Range("J138").Value = 5255.375
$ws.Range("K138").Value = 6696.999899999999
$ws.Range("L138").Value = 15766.125
$ws.Range("M138").Value = -1556.999899999999
$ws.Range("N138").Value = -26046.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1032.6666
$ws.Range("J5").Value = 1524
$ws.Range("L5").Value = 1524
$ws.Range("N5").Value = -1748
$ws.Range("H32").Value = 6494.1445
$ws.Range("I32").Value = 6264.257
$ws.Range("K32").Value = 6264.257
$ws.Range("M32").Value = -5977.257
$ws.Range("H45").Value = 344833.34
$ws.Range("I45").Value = 1015000
$ws.Range("J45").Value = 9750
$ws.Range("K45").Value = 1015000
$ws.Range("L45").Value = 9750
$ws.Range("M45").Value = -1014623
$ws.Range("N45").Value = -10504
$ws.Range("H74").Value = 6789.905
$ws.Range("I74").Value = 17399.5
$ws.Range("J74").Value = 2546.0667
$ws.Range("K74").Value = 17399.5
$ws.Range("L74").Value = 2546.0667
$ws.Range("M74").Value = -16525.5
$ws.Range("N74").Value = -4294.066699999999
$ws.Range("H77").Value = 6789.905
$ws.Range("I77").Value = 17399.5
$ws.Range("J77").Value = 2546.0667
$ws.Range("K77").Value = 86997.5
$ws.Range("L77").Value = 12730.3335
$ws.Range("M77").Value = -82629.5
$ws.Range("N77").Value = -21466.3335
$ws.Range("H132").Value = 3023.532
$ws.Range("I132").Value = 2040.0938
$ws.Range("J132").Value = 5121.533
$ws.Range("K132").Value = 6120.2814
$ws.Range("L132").Value = 15364.599
$ws.Range("M132").Value = -3590.2814
$ws.Range("N132").Value = -20424.599

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1032.6666
$ws.Range("J4").Value = 1524
$ws.Range("L4").Value = 1524
$ws.Range("N4").Value = -1754
$ws.Range("H134").Value = 10323.277
$ws.Range("I134").Value = 11781.786
$ws.Range("J134").Value = 5218.5
$ws.Range("K134").Value = 35345.358
$ws.Range("L134").Value = 15655.5
$ws.Range("M134").Value = -32810.358
$ws.Range("N134").Value = -20725.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6363.6
$ws.Range("I31").Value = 8035.385
$ws.Range("J31").Value = 4864.759
$ws.Range("K31").Value = 8035.385
$ws.Range("L31").Value = 4864.759
$ws.Range("M31").Value = -7740.385
$ws.Range("N31").Value = -5454.759
$ws.Range("H34").Value = 6363.6
$ws.Range("I34").Value = 8035.385
$ws.Range("J34").Value = 4864.759
$ws.Range("K34").Value = 8035.385
$ws.Range("L34").Value = 4864.759
$ws.Range("M34").Value = -7833.385
$ws.Range("N34").Value = -5268.759
$ws.Range("H99").Value = 253822.25
$ws.Range("I99").Value = 502344.5
$ws.Range("J99").Value = 5300
$ws.Range("K99").Value = 502344.5
$ws.Range("L99").Value = 5300
$ws.Range("M99").Value = -500846.5
$ws.Range("N99").Value = -8296
$ws.Range("H122").Value = 12500.637
$ws.Range("I122").Value = 21276.334
$ws.Range("K122").Value = 63829.00199999999
$ws.Range("M122").Value = -61379.00199999999
$ws.Range("H126").Value = 253822.25
$ws.Range("I126").Value = 502344.5
$ws.Range("J126").Value = 5300
$ws.Range("K126").Value = 1507033.5
$ws.Range("L126").Value = 15900
$ws.Range("M126").Value = -1504563.5
$ws.Range("N126").Value = -20840
$ws.Range("H134").Value = 4765.9688
$ws.Range("I134").Value = 4986.483
$ws.Range("K134").Value = 14959.449
$ws.Range("M134").Value = -12424.449
$ws.Range("H141").Value = 188333.72
$ws.Range("J141").Value = 202533.77
$ws.Range("L141").Value = 202533.77
$ws.Range("N141").Value = -212893.77

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 16000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 16000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 48000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -49118
$ws.Range("H80").Value = 147013.8
$ws.Range("J80").Value = 96889.78
$ws.Range("L80").Value = 290669.34
$ws.Range("N80").Value = -292541.34
$ws.Range("H83").Value = 147013.8
$ws.Range("J83").Value = 96889.78
$ws.Range("L83").Value = 872008.02
$ws.Range("N83").Value = -881368.02
$ws.Range("H114").Value = 2524.9092
$ws.Range("J114").Value = 3996.25
$ws.Range("L114").Value = 11988.75
$ws.Range("N114").Value = -18496.75
$ws.Range("H131").Value = 10640969
$ws.Range("I131").Value = 76929790
$ws.Range("J131").Value = 2021.2839
$ws.Range("K131").Value = 230789370
$ws.Range("L131").Value = 6063.851699999999
$ws.Range("M131").Value = -230784330
$ws.Range("N131").Value = -16143.8517

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9868.5
$ws.Range("I132").Value = 9868.5
$ws.Range("K132").Value = 29605.5
$ws.Range("M132").Value = -27075.5
$ws.Range("H139").Value = 46084.223
$ws.Range("J139").Value = 46084.223
$ws.Range("L139").Value = 46084.223
$ws.Range("N139").Value = -56364.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 18588.676
$ws.Range("I40").Value = 20957.826
$ws.Range("J40").Value = 13635
$ws.Range("K40").Value = 20957.826
$ws.Range("L40").Value = 13635
$ws.Range("M40").Value = -20821.826
$ws.Range("N40").Value = -13907
$ws.Range("H61").Value = 3082.0605
$ws.Range("I61").Value = 1376.35
$ws.Range("K61").Value = 1376.35
$ws.Range("M61").Value = -1174.35
$ws.Range("H113").Value = 3082.0605
$ws.Range("I113").Value = 1376.35
$ws.Range("K113").Value = 1376.35
$ws.Range("M113").Value = 793.6500000000001
$ws.Range("H132").Value = 650899.4
$ws.Range("I132").Value = 1243036.1
$ws.Range("J132").Value = 4932.091
$ws.Range("K132").Value = 3729108.3
$ws.Range("L132").Value = 14796.273
$ws.Range("M132").Value = -3726578.3
$ws.Range("N132").Value = -19856.273
$ws.Range("H136").Value = 7070.9414
$ws.Range("J136").Value = 6813.8
$ws.Range("L136").Value = 20441.4
$ws.Range("N136").Value = -25541.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H113").Value = 1003.06665
$ws.Range("I113").Value = 538.57574
$ws.Range("J113").Value = 2280.4167
$ws.Range("K113").Value = 1615.72722
$ws.Range("L113").Value = 6841.250100000001
$ws.Range("M113").Value = 554.27278
$ws.Range("N113").Value = -11181.2501
$ws.Range("H122").Value = 15792.775
$ws.Range("I122").Value = 1828
$ws.Range("J122").Value = 39067.4
$ws.Range("K122").Value = 5484
$ws.Range("L122").Value = 117202.2
$ws.Range("M122").Value = -3034
$ws.Range("N122").Value = -122102.2
$ws.Range("H132").Value = 7620.2905
$ws.Range("I132").Value = 8541.022000000001
$ws.Range("J132").Value = 5183.0586
$ws.Range("K132").Value = 25623.066
$ws.Range("L132").Value = 15549.1758
$ws.Range("M132").Value = -23093.066
$ws.Range("N132").Value = -20609.1758
$ws.Range("H136").Value = 679792.75
$ws.Range("I136").Value = 913425.6
$ws.Range("J136").Value = 17832.834
$ws.Range("K136").Value = 2740276.8
$ws.Range("L136").Value = 53498.50199999999
$ws.Range("M136").Value = -2737726.8
$ws.Range("N136").Value = -58598.50199999999
